$d = $word.ActiveDocument

# 1) Replace the trailing period with an ellipsis character in the
#    "Implementando Lista Persistente de artículos." run.
$d.Content.Find.Execute("Lista Persistente de artículos.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Lista Persistente de artículos…", 2)

# 2) Merge that paragraph with the following bookmark-only ("_GoBack")
#    paragraph by deleting the paragraph mark that separates them, so the
#    bookmark ends up living in the same paragraph as the text above.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Lista Persistente de artículos…*") {
        $target = $p
        break
    }
}
$mergeRange = $d.Range($target.Range.End - 1, $target.Range.End)
$mergeRange.Delete()

# 3) Insert a new empty paragraph right after the merged paragraph so the
#    overall paragraph count stays the same as before the merge.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Lista Persistente de artículos…*") {
        $target = $p
        break
    }
}
$target.Range.InsertParagraphAfter()
